$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H52").Value = 3169.5217
$ws.Range("J52").Value = 3222.682
$ws.Range("L52").Value = 9668.045999999998
$ws.Range("N52").Value = -9988.045999999998
$ws.Range("H135").Value = 4303.2
$ws.Range("I135").Value = 3670.2222
$ws.Range("J135").Value = 10000
$ws.Range("K135").Value = 33031.99980000001
$ws.Range("L135").Value = 90000
$ws.Range("M135").Value = -30496.99980000001
$ws.Range("N135").Value = -95070
$ws.Range("H137").Value = 2821.4468
$ws.Range("I137").Value = 2580.3823
$ws.Range("J137").Value = 3451.923
$ws.Range("K137").Value = 7741.146900000001
$ws.Range("L137").Value = 10355.769
$ws.Range("M137").Value = -5191.146900000001
$ws.Range("N137").Value = -15455.769
$ws.Range("H138").Value = 1874.7931
$ws.Range("J138").Value = 1878.9697
$ws.Range("L138").Value = 5636.909100000001
$ws.Range("N138").Value = -15916.9091

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H28").Value = 10742.75
$ws.Range("I28").Value = 4323.6665
$ws.Range("K28").Value = 4323.6665
$ws.Range("M28").Value = -4131.6665
$ws.Range("H45").Value = 2174.1177
$ws.Range("I45").Value = 1475.25
$ws.Range("J45").Value = 3851.4
$ws.Range("K45").Value = 1475.25
$ws.Range("L45").Value = 3851.4
$ws.Range("M45").Value = -1098.25
$ws.Range("N45").Value = -4605.4
$ws.Range("H74").Value = 1594.5714
$ws.Range("I74").Value = 1616.0588
$ws.Range("J74").Value = 1503.25
$ws.Range("K74").Value = 1616.0588
$ws.Range("L74").Value = 1503.25
$ws.Range("M74").Value = -742.0588
$ws.Range("N74").Value = -3251.25
$ws.Range("H77").Value = 1594.5714
$ws.Range("I77").Value = 1616.0588
$ws.Range("J77").Value = 1503.25
$ws.Range("K77").Value = 8080.294
$ws.Range("L77").Value = 7516.25
$ws.Range("M77").Value = -3712.294
$ws.Range("N77").Value = -16252.25
$ws.Range("H99").Value = 10742.75
$ws.Range("I99").Value = 4323.6665
$ws.Range("K99").Value = 4323.6665
$ws.Range("M99").Value = -1328.6665

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 2710.2693
$ws.Range("I86").Value = 2836.8125
$ws.Range("J86").Value = 2507.8
$ws.Range("K86").Value = 2836.8125
$ws.Range("L86").Value = 2507.8
$ws.Range("M86").Value = -1713.8125
$ws.Range("N86").Value = -4753.8
$ws.Range("H89").Value = 2710.2693
$ws.Range("I89").Value = 2836.8125
$ws.Range("J89").Value = 2507.8
$ws.Range("K89").Value = 14184.0625
$ws.Range("L89").Value = 12539
$ws.Range("M89").Value = -8568.0625
$ws.Range("N89").Value = -23771
$ws.Range("H99").Value = 1904.4445
$ws.Range("I99").Value = 1752.5
$ws.Range("J99").Value = 1968.421
$ws.Range("K99").Value = 1752.5
$ws.Range("L99").Value = 1968.421
$ws.Range("M99").Value = -254.5
$ws.Range("N99").Value = -4964.421
$ws.Range("H126").Value = 1904.4445
$ws.Range("I126").Value = 1752.5
$ws.Range("J126").Value = 1968.421
$ws.Range("K126").Value = 5257.5
$ws.Range("L126").Value = 5905.263
$ws.Range("M126").Value = -2787.5
$ws.Range("N126").Value = -10845.263
$ws.Range("H132").Value = 5954123.5
$ws.Range("I132").Value = 1144.7
$ws.Range("K132").Value = 3434.1
$ws.Range("M132").Value = -904.1000000000004

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H118").Value = 2838.9387
$ws.Range("I118").Value = 1164.5
$ws.Range("J118").Value = 2987.7778
$ws.Range("K118").Value = 3493.5
$ws.Range("L118").Value = 8963.3334
$ws.Range("M118").Value = -2250.5
$ws.Range("N118").Value = -11449.3334
$ws.Range("H137").Value = 16357.75
$ws.Range("J137").Value = 4311
$ws.Range("L137").Value = 12933
$ws.Range("N137").Value = -23133
$ws.Range("H140").Value = 1422.8462
$ws.Range("I140").Value = 956.6957
$ws.Range("J140").Value = 4996.6665
$ws.Range("K140").Value = 2870.0871
$ws.Range("L140").Value = 14989.9995
$ws.Range("M140").Value = 2309.9129
$ws.Range("N140").Value = -25349.9995

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 4362.3516
$ws.Range("I122").Value = 1565.1538
$ws.Range("J122").Value = 5877.5
$ws.Range("K122").Value = 4695.4614
$ws.Range("L122").Value = 17632.5
$ws.Range("M122").Value = -2245.4614
$ws.Range("N122").Value = -22532.5
$ws.Range("H126").Value = 3602.8
$ws.Range("I126").Value = 3202.8
$ws.Range("K126").Value = 9608.400000000001
$ws.Range("M126").Value = -7138.400000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4577.722
$ws.Range("I40").Value = 1774.75
$ws.Range("K40").Value = 1774.75
$ws.Range("M40").Value = -1638.75
$ws.Range("H132").Value = 4915.615
$ws.Range("I132").Value = 4020.1667
$ws.Range("J132").Value = 5683.143
$ws.Range("K132").Value = 12060.5001
$ws.Range("L132").Value = 17049.429
$ws.Range("M132").Value = -9530.500100000001
$ws.Range("N132").Value = -22109.429
$ws.Range("H136").Value = 2249.682
$ws.Range("I136").Value = 2426.2
$ws.Range("J136").Value = 1871.4286
$ws.Range("K136").Value = 7278.599999999999
$ws.Range("L136").Value = 5614.2858
$ws.Range("M136").Value = -4728.599999999999
$ws.Range("N136").Value = -10714.2858

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 4176.533
$ws.Range("I81").Value = 4365
$ws.Range("J81").Value = 3799.6
$ws.Range("K81").Value = 8730
$ws.Range("L81").Value = 7599.2
$ws.Range("M81").Value = -7669
$ws.Range("N81").Value = -9721.200000000001
$ws.Range("H84").Value = 4176.533
$ws.Range("I84").Value = 4365
$ws.Range("J84").Value = 3799.6
$ws.Range("K84").Value = 43650
$ws.Range("L84").Value = 37996
$ws.Range("M84").Value = -38346
$ws.Range("N84").Value = -48604
$ws.Range("H122").Value = 3387.9
$ws.Range("I122").Value = 964.7143
$ws.Range("K122").Value = 2894.1429
$ws.Range("M122").Value = -444.1428999999998
$ws.Range("H126").Value = 1682.6428
$ws.Range("I126").Value = 1507.4445
$ws.Range("K126").Value = 4522.333500000001
$ws.Range("M126").Value = -2052.333500000001
$ws.Range("H136").Value = 2836.1714
$ws.Range("I136").Value = 2591.55
$ws.Range("J136").Value = 3162.3333
$ws.Range("K136").Value = 7774.650000000001
$ws.Range("L136").Value = 9486.999899999999
$ws.Range("M136").Value = -5224.650000000001
$ws.Range("N136").Value = -14586.9999
